# Commit: modifying the items and creating new Class
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 4 quantity (C4): 35 -> 40
$ws.Range("C4").Value = 40

# Add new item row 6: 12313 / Milk / 30
$ws.Range("A6").Value = 12313
$ws.Range("B6").Value = "Milk"
$ws.Range("C6").Value = 30

# Add new item row 7: 11111 / Smokic / 10
$ws.Range("A7").Value = 11111
$ws.Range("B7").Value = "Smokic"
$ws.Range("C7").Value = 10

# Move the active selection to C4, matching the saved view state
$ws.Range("C4").Select() | Out-Null
